# Mise à jour de l'application
# Appends the 12 new GPS data rows (match "N3 J9 VS Beaucaire", 2025-11-22)
# to the bottom of the Feuil1 sheet, then moves the view/selection down
# to where the new data was entered (matching the author's final cursor
# position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$matchLabel = "N3 J9 VS Beaucaire"
$matchDate  = 45983
$periodType = "Global"
$md         = "M"

# row: @(RowNum, Player, Poste, TempsJoue,
#         Distance, DistanceHID, D0_15, D15_20, D20_25, D25_30, D30plus,
#         Sprints, VitMoy, VitMax, AccelMax, Accel3, Accel4, Decel3, Decel4)
$rows = @(
    @(898, "Kamal Bafounta",   "center midfield", "00:09:48", 1.24,  0.18, 1.07, 0.16, 0.02, 0,    0,    0,  7.39, 22.77, 3.74, 6,  0,  5,  1),
    @(899, "Mattheo Haon",     "right back",       "01:35:45", 11.4,  2.06, 9.31, 1.31, 0.56, 0.19, 0.03, 14, 6.98, 31.82, 5.03, 36, 6,  38, 14),
    @(900, "Emmanuel Valey",   "left forward",     "01:07:40", 8.4,   2.36, 6.01, 1.45, 0.75, 0.18, 0.01, 18, 7.4,  31.86, 4.6,  47, 8,  37, 14),
    @(901, "Yoan Zouma",       "center back",      "01:35:38", 9.22,  0.9,  8.31, 0.67, 0.22, 0.02, 0,    2,  5.63, 26.67, 3.96, 19, 0,  23, 1),
    @(902, "Amir Etien",       "right forward",    "01:25:26", 8.12,  1.92, 6.18, 1.02, 0.65, 0.2,  0.06, 15, 5.62, 33.54, 5.32, 40, 17, 24, 14),
    @(903, "Naim Ighbane",     "center back",      "01:34:34", 9.86,  1.29, 8.55, 0.9,  0.32, 0.07, 0.01, 5,  5.96, 31.21, 4.39, 14, 1,  14, 4),
    @(904, "Maé Clavel",       "left back",        "00:16:20", 2.07,  0.5,  1.56, 0.31, 0.11, 0.08, 0,    3,  7.61, 29.38, 4.1,  8,  1,  4,  3),
    @(905, "Jeremie Laurent",  "left forward",     "00:27:34", 3.41,  1.01, 2.38, 0.52, 0.34, 0.16, 0,    10, 7.44, 30.13, 4.1,  14, 2,  17, 6),
    @(906, "Malik Boussaid",   "right back",       "01:18:54", 9.54,  1.92, 7.6,  1.27, 0.52, 0.15, 0,    8,  7.17, 29.61, 4.5,  33, 6,  36, 6),
    @(907, "Yoann Martelat",   "center midfield",  "01:36:09", 12.03, 2.77, 9.23, 2.27, 0.45, 0.08, 0,    6,  7.36, 27.96, 5.79, 25, 6,  29, 5),
    @(908, "Naim Dhib",        "center midfield",  "01:36:01", 10.51, 2.2,  8.28, 1.56, 0.56, 0.11, 0,    8,  6.51, 29.58, 5.03, 61, 6,  40, 8),
    @(909, "Sofiane Belle",    "left forward",     "01:19:57", 7.81,  1.63, 6.15, 0.97, 0.45, 0.23, 0.01, 14, 5.77, 30.63, 4.63, 18, 2,  39, 9)
)

# Two passes so new shared-string entries land in the same order the
# original author's save produced them: the "Temps joué" (G) values for
# each row first, then the repeated match label (A) last (it only mints
# one new shared-string entry, on its first use).
foreach ($row in $rows) {
    $r = $row[0]

    # --- Column B: date, copy the date-number-format style from the row above ---
    $ws.Cells.Item($r - 1, 2).Copy()
    $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    $ws.Cells.Item($r, 2).Value = $matchDate

    # --- Column C: period/type ---
    $ws.Cells.Item($r, 3).Value = $periodType

    # --- Column D: MD, copy the centered style from the row above ---
    $ws.Cells.Item($r - 1, 4).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
    $ws.Cells.Item($r, 4).Value = $md

    # --- Column E: player name ---
    $ws.Cells.Item($r, 5).Value = $row[1]

    # --- Column F: position ---
    $ws.Cells.Item($r, 6).Value = $row[2]

    # --- Column G: temps joué ---
    $ws.Cells.Item($r, 7).Value = $row[3]

    # --- Columns H..V: numeric GPS metrics ---
    for ($i = 0; $i -lt 15; $i++) {
        $ws.Cells.Item($r, 8 + $i).Value = $row[4 + $i]
    }
}

foreach ($row in $rows) {
    $r = $row[0]
    # --- Column A: match/session label (plain string, same style as neighbour) ---
    $ws.Cells.Item($r, 1).Value = $matchLabel
}

$excel.CutCopyMode = 0

# Move the view/selection to reflect where data entry ended.
$ws.Range("D914").Select()
